$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("A1")
$v = $r.get_Value()
Write-Host "GET_VALUE:" $v
$v2 = $r.Value()
Write-Host "VALUE():" $v2
